$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 63: 'Summoning for Dummies' / 'Archaeoskin Codex'
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# Row 66: 'Summoning the Courage to Be Different (L)' / 'Archaeoskin Codex'
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# Row 76: 'Warding Off Temptation' / 'Enchanted Hardsilver Ink'
$ws.Range("H76").Value = 3000
$ws.Range("I76").Value = 3000
$ws.Range("K76").Value = 3000
$ws.Range("M76").Value = -2685

# Row 79: 'The Garden of Arcane Delights (L)' / 'Enchanted Hardsilver Ink'
$ws.Range("H79").Value = 3000
$ws.Range("I79").Value = 3000
$ws.Range("K79").Value = 3000
$ws.Range("M79").Value = -1908

# Row 106: 'Making Your Mark' / 'Enchanted Palladium Ink'
$ws.Range("H106").Value = 720
$ws.Range("I106").Value = 720
$ws.Range("K106").Value = 720
$ws.Range("M106").Value = -89

# Row 137: 'Cutting Edge of Culinary Quality' / 'Magnesia Whetstone'
$ws.Range("H137").Value = 1377.4117
$ws.Range("I137").Value = 1276
$ws.Range("K137").Value = 3828
$ws.Range("M137").Value = -1278

# Row 138: 'All-night Crafting' / "Cunning Craftsman's Tisane"
$ws.Range("H138").Value = 12597.464
$ws.Range("I138").Value = 9747.5
$ws.Range("J138").Value = 12743.615
$ws.Range("K138").Value = 29242.5
$ws.Range("L138").Value = 38230.845
$ws.Range("M138").Value = -24102.5
$ws.Range("N138").Value = -48510.845

# Row 141: 'Remedy for Reason' / 'Grade 1 Gemdraught of Mind'
$ws.Range("H141").Value = 1098.75
$ws.Range("I141").Value = 1098.75
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 3296.25
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 1883.75
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 2: "Ain't Got No Ingots" / 'Bronze Ingot'
$ws.Range("H2").Value = 1142.7142
$ws.Range("I2").Value = 1142.7142
$ws.Range("K2").Value = 1142.7142
$ws.Range("M2").Value = -1029.7142

# Row 32: 'Ingot We Trust' / 'Steel Ingot'
$ws.Range("H32").Value = 9590.868
$ws.Range("I32").Value = 9174.433000000001
$ws.Range("K32").Value = 9174.433000000001
$ws.Range("M32").Value = -8887.433000000001

# Row 61: 'Dealing with the Tough Stuff' / 'Cobalt Ingot'
$ws.Range("H61").Value = 3450.8572
$ws.Range("J61").Value = 4799.5713
$ws.Range("L61").Value = 4799.5713
$ws.Range("N61").Value = -5223.5713

# Row 74: 'As the Bolt Flies' / 'Titanium Nugget'
$ws.Range("H74").Value = 2799.5
$ws.Range("I74").Value = 2799.5
$ws.Range("K74").Value = 2799.5
$ws.Range("M74").Value = -1925.5

# Row 77: 'Heavy Metal Banned (L)' / 'Titanium Nugget'
$ws.Range("H77").Value = 2799.5
$ws.Range("I77").Value = 2799.5
$ws.Range("K77").Value = 13997.5
$ws.Range("M77").Value = -9629.5

# Row 116: 'No Scope' / 'Titanbronze Ingot'
$ws.Range("H116").Value = 1142.7142
$ws.Range("I116").Value = 1142.7142
$ws.Range("K116").Value = 1142.7142
$ws.Range("M116").Value = 1151.2858

# Row 132: "Don't Bore Me, Ore Me" / 'Mountain Chromite Ingot'
$ws.Range("H132").Value = 1359.8529
$ws.Range("I132").Value = 1085.16
$ws.Range("J132").Value = 2122.889
$ws.Range("K132").Value = 3255.48
$ws.Range("L132").Value = 6368.667
$ws.Range("M132").Value = -725.4800000000005
$ws.Range("N132").Value = -11428.667

# Row 136: 'Metal with Mettle' / 'Cobalt Tungsten Ingot'
$ws.Range("H136").Value = 3450.8572
$ws.Range("J136").Value = 4799.5713
$ws.Range("L136").Value = 14398.7139
$ws.Range("N136").Value = -19498.7139

$ws = $wb.Worksheets.Item("BSM")
# Row 3: 'Hells Bells' / 'Bronze Ingot'
$ws.Range("H3").Value = 1142.7142
$ws.Range("I3").Value = 1142.7142
$ws.Range("K3").Value = 1142.7142
$ws.Range("M3").Value = -1028.7142

# Row 22: 'Riveting Run' / 'Iron Rivets'
$ws.Range("H22").Value = 537.4
$ws.Range("I22").Value = 586.1111
$ws.Range("K22").Value = 586.1111
$ws.Range("M22").Value = -413.1111

# Row 29: 'Powderpost Derby' / "Initiate's Saw"
$ws.Range("H29").Value = 459.33334
$ws.Range("I29").Value = 459.33334
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 459.33334
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -170.33334
$ws.Range("N29").ClearContents()

# Row 99: 'Meddle in Metal' / 'Oroshigane Ingot'
$ws.Range("H99").Value = 1403.3334
$ws.Range("I99").Value = 210
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 210
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 1288
$ws.Range("N99").Value = -4996

# Row 105: 'Ingot to Wing It' / 'Molybdenum Ingot'
$ws.Range("H105").Value = 4926.231
$ws.Range("I105").Value = 3934.5715
$ws.Range("K105").Value = 3934.5715
$ws.Range("M105").Value = -2187.5715

# Row 107: 'The Gold Experience' / 'Deepgold Nugget'
$ws.Range("H107").Value = 1098.9231
$ws.Range("I107").Value = 1065.5
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1065.5
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 854.5
$ws.Range("N107").Value = -5340

# Row 134: 'Ruthenium Supremium' / 'Ruthenium Ingot'
$ws.Range("H134").Value = 1681.6364
$ws.Range("I134").Value = 1512.25
$ws.Range("K134").Value = 4536.75
$ws.Range("M134").Value = -2001.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31: 'Wall Not Found' / 'Walnut Lumber'
$ws.Range("H31").Value = 3272.6875
$ws.Range("I31").Value = 3219
$ws.Range("J31").Value = 3648.5
$ws.Range("K31").Value = 3219
$ws.Range("L31").Value = 3648.5
$ws.Range("M31").Value = -2924
$ws.Range("N31").Value = -4238.5

# Row 34: 'Armoires of the Rich and Famous' / 'Walnut Lumber'
$ws.Range("H34").Value = 3272.6875
$ws.Range("I34").Value = 3219
$ws.Range("J34").Value = 3648.5
$ws.Range("K34").Value = 3219
$ws.Range("L34").Value = 3648.5
$ws.Range("M34").Value = -3017
$ws.Range("N34").Value = -4052.5

# Row 86: 'Birch, Please' / 'Birch Lumber'
$ws.Range("H86").Value = 8715857
$ws.Range("I86").Value = 9959554
$ws.Range("K86").Value = 9959554
$ws.Range("M86").Value = -9958431

# Row 89: 'Built This City on Blocks and Soul (L)' / 'Birch Lumber'
$ws.Range("H89").Value = 8715857
$ws.Range("I89").Value = 9959554
$ws.Range("K89").Value = 49797770
$ws.Range("M89").Value = -49792154

# Row 107: 'Built to Last' / 'White Oak Lumber'
$ws.Range("H107").Value = 853.14813
$ws.Range("I107").Value = 523.34784
$ws.Range("J107").Value = 2749.5
$ws.Range("K107").Value = 523.34784
$ws.Range("L107").Value = 2749.5
$ws.Range("M107").Value = 1396.65216
$ws.Range("N107").Value = -6589.5

$ws = $wb.Worksheets.Item("CUL")
# Row 2: 'Pork Is a Salty Food' / 'Table Salt'
$ws.Range("H2").Value = 43.384617
$ws.Range("I2").Value = 20.333334
$ws.Range("K2").Value = 122.000004
$ws.Range("M2").Value = -9.000004000000004

$ws = $wb.Worksheets.Item("GSM")
# Row 43: 'Get the Green Stuff' / 'Malachite Earrings'
$ws.Range("H43").Value = 17929.25
$ws.Range("I43").Value = 7239
$ws.Range("K43").Value = 7239
$ws.Range("M43").Value = -7088

# Row 63: 'Not on My Table' / 'Mythrite Earrings of Healing'
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# Row 66: "Heinz's Dilemma (L)" / 'Mythrite Earrings of Healing'
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# Row 70: 'Sky Is the Limit' / 'Mythrite Ingot'
$ws.Range("H70").Value = 25002732
$ws.Range("I70").Value = 25002732
$ws.Range("K70").Value = 25002732
$ws.Range("M70").Value = -25002462

# Row 73: 'Hulls of Broken Dreams (L)' / 'Mythrite Ingot'
$ws.Range("H73").Value = 25002732
$ws.Range("I73").Value = 25002732
$ws.Range("K73").Value = 25002732
$ws.Range("M73").Value = -25001796

# Row 102: 'Put the Metal to the Peddle' / 'Durium Ingot'
$ws.Range("H102").Value = 1240.1111
$ws.Range("I102").Value = 995.2857
$ws.Range("K102").Value = 995.2857
$ws.Range("M102").Value = 626.7143

# Row 113: 'Copious Crystal Cannons' / 'Manasilver Nugget'
$ws.Range("H113").Value = 1231.5
$ws.Range("I113").Value = 1231.5
$ws.Range("K113").Value = 1231.5
$ws.Range("M113").Value = 938.5

# Row 122: 'Awarding Academic Excellence' / 'Ametrine'
$ws.Range("H122").Value = 5752.0586
$ws.Range("I122").Value = 5799.0625
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 17397.1875
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -14947.1875
$ws.Range("N122").Value = -19900

# Row 126: 'Gold Rush Order' / 'Phrygian Gold Ingot'
$ws.Range("H126").Value = 250003400
$ws.Range("J126").Value = 5750
$ws.Range("L126").Value = 17250
$ws.Range("N126").Value = -22190

# Row 132: 'On Board for Lar' / 'Lar Ingot'
$ws.Range("H132").Value = 2498.125
$ws.Range("I132").Value = 2121.6667
$ws.Range("K132").Value = 6365.000100000001
$ws.Range("M132").Value = -3835.000100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 22: 'Skin off Their Backs' / 'Aldgoat Leather'
$ws.Range("H22").Value = 3195
$ws.Range("J22").Value = 4200
$ws.Range("L22").Value = 4200
$ws.Range("N22").Value = -4790

# Row 27: 'Fire and Hide' / 'Aldgoat Leather'
$ws.Range("H27").Value = 3195
$ws.Range("J27").Value = 4200
$ws.Range("L27").Value = 4200
$ws.Range("N27").Value = -4414

# Row 40: 'Best Served Toad' / 'Toad Leather'
$ws.Range("H40").Value = 12430.3125
$ws.Range("I40").Value = 12420.429
$ws.Range("J40").Value = 12499.5
$ws.Range("K40").Value = 12420.429
$ws.Range("L40").Value = 12499.5
$ws.Range("M40").Value = -12284.429
$ws.Range("N40").Value = -12771.5

# Row 55: "It's Not a Job, It's a Calling" / 'Peiste Leather'
$ws.Range("H55").Value = 903.7222
$ws.Range("I55").Value = 362.85715
$ws.Range("J55").Value = 1247.909
$ws.Range("K55").Value = 362.85715
$ws.Range("L55").Value = 1247.909
$ws.Range("M55").Value = -189.85715
$ws.Range("N55").Value = -1593.909

# Row 61: 'Spelling Me Softly' / 'Raptor Leather'
$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

# Row 113: 'Peace in Rest' / 'Atrociraptor Leather'
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 64: 'Ribbon of Remembrance' / 'Rainbow Ribbon of Healing'
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67: 'The Road Was a Ribbon of Moonlight (L)' / 'Rainbow Ribbon of Healing'
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 113: 'A Tender Table' / 'Pixie Floss'
$ws.Range("H113").Value = 857.8823
$ws.Range("I113").Value = 630.36365
$ws.Range("K113").Value = 1891.09095
$ws.Range("M113").Value = 278.90905

# Row 122: 'Heavy Armoire' / 'Dark Hempen Cloth'
$ws.Range("H122").Value = 975
$ws.Range("J122").Value = 1199
$ws.Range("L122").Value = 3597
$ws.Range("N122").Value = -8497

# Row 126: 'A Polished Purchase' / 'Snow Linen'
$ws.Range("H126").Value = 2999.3333
$ws.Range("I126").Value = 2333
$ws.Range("K126").Value = 6999
$ws.Range("M126").Value = -4529

# Row 136: 'Weaving the Envelope' / 'Sarcenet Cloth'
$ws.Range("H136").Value = 29346.37
$ws.Range("I136").Value = 33598.875
$ws.Range("K136").Value = 100796.625
$ws.Range("M136").Value = -98246.625
